# Swap the order of "Recorded By" names in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
